$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1750
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 2500
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 2500
$ws.Range("M25").Value = -598
$ws.Range("N25").Value = -3304

$ws.Range("H74").Value = 16130234
$ws.Range("I74").Value = 18519658
$ws.Range("J74").Value = 1617
$ws.Range("K74").Value = 18519658
$ws.Range("L74").Value = 1617
$ws.Range("M74").Value = -18518784
$ws.Range("N74").Value = -3365

$ws.Range("H77").Value = 16130234
$ws.Range("I77").Value = 18519658
$ws.Range("J77").Value = 1617
$ws.Range("K77").Value = 92598290
$ws.Range("L77").Value = 8085
$ws.Range("M77").Value = -92593922
$ws.Range("N77").Value = -16821

$ws.Range("H92").Value = 50550
$ws.Range("J92").Value = 50550
$ws.Range("L92").Value = 50550
$ws.Range("N92").Value = -55542

$ws.Range("H97").Value = 6091.7896
$ws.Range("J97").Value = 1403.4286
$ws.Range("L97").Value = 1403.4286
$ws.Range("N97").Value = -2395.4286

$ws.Range("H102").Value = 2662.4443
$ws.Range("I102").Value = 2532.8
$ws.Range("K102").Value = 2532.8
$ws.Range("M102").Value = -910.8000000000002

$ws.Range("H132").Value = 6412046.5
$ws.Range("I132").Value = 8930074
$ws.Range("J132").Value = 2520.9092
$ws.Range("K132").Value = 26790222
$ws.Range("L132").Value = 7562.7276
$ws.Range("M132").Value = -26787692
$ws.Range("N132").Value = -12622.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1748.8572
$ws.Range("I20").Value = 1789.4762
$ws.Range("J20").Value = 1627
$ws.Range("K20").Value = 1789.4762
$ws.Range("L20").Value = 1627
$ws.Range("M20").Value = -1542.4762
$ws.Range("N20").Value = -2121

$ws.Range("H86").Value = 29414560
$ws.Range("I86").Value = 2636.818
$ws.Range("K86").Value = 2636.818
$ws.Range("M86").Value = -1513.818

$ws.Range("H89").Value = 29414560
$ws.Range("I89").Value = 2636.818
$ws.Range("K89").Value = 13184.09
$ws.Range("M89").Value = -7568.09

$ws.Range("H94").Value = 611.2564
$ws.Range("I94").Value = 573.2963
$ws.Range("J94").Value = 696.6667
$ws.Range("K94").Value = 573.2963
$ws.Range("L94").Value = 696.6667
$ws.Range("M94").Value = -122.2963
$ws.Range("N94").Value = -1598.6667

$ws.Range("H105").Value = 3065.7346
$ws.Range("I105").Value = 1620.92
$ws.Range("J105").Value = 4570.75
$ws.Range("K105").Value = 1620.92
$ws.Range("L105").Value = 4570.75
$ws.Range("M105").Value = 126.0799999999999
$ws.Range("N105").Value = -8064.75

$ws.Range("H107").Value = 799.1818
$ws.Range("I107").Value = 754.55554
$ws.Range("K107").Value = 754.55554
$ws.Range("M107").Value = 1165.44446

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8551574
$ws.Range("I31").Value = 4712.3438
$ws.Range("J31").Value = 47622940
$ws.Range("K31").Value = 4712.3438
$ws.Range("L31").Value = 47622940
$ws.Range("M31").Value = -4417.3438
$ws.Range("N31").Value = -47623530

$ws.Range("H34").Value = 8551574
$ws.Range("I34").Value = 4712.3438
$ws.Range("J34").Value = 47622940
$ws.Range("K34").Value = 4712.3438
$ws.Range("L34").Value = 47622940
$ws.Range("M34").Value = -4510.3438
$ws.Range("N34").Value = -47623344

$ws.Range("H58").Value = 1695.7142
$ws.Range("I58").Value = 781.1739
$ws.Range("J58").Value = 3448.5833
$ws.Range("K58").Value = 781.1739
$ws.Range("L58").Value = 3448.5833
$ws.Range("M58").Value = -578.1739
$ws.Range("N58").Value = -3854.5833

$ws.Range("H99").Value = 1696.9131
$ws.Range("I99").Value = 1618.2778
$ws.Range("J99").Value = 1980
$ws.Range("K99").Value = 1618.2778
$ws.Range("L99").Value = 1980
$ws.Range("M99").Value = -120.2778000000001
$ws.Range("N99").Value = -4976

$ws.Range("H126").Value = 1696.9131
$ws.Range("I126").Value = 1618.2778
$ws.Range("J126").Value = 1980
$ws.Range("K126").Value = 4854.8334
$ws.Range("L126").Value = 5940
$ws.Range("M126").Value = -2384.8334
$ws.Range("N126").Value = -10880

$ws.Range("H132").Value = 17859474
$ws.Range("I132").Value = 22729290
$ws.Range("K132").Value = 68187870
$ws.Range("M132").Value = -68185340

$ws.Range("H136").Value = 1695.7142
$ws.Range("I136").Value = 781.1739
$ws.Range("J136").Value = 3448.5833
$ws.Range("K136").Value = 2343.5217
$ws.Range("L136").Value = 10345.7499
$ws.Range("M136").Value = 206.4782999999998
$ws.Range("N136").Value = -15445.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 156.2
$ws.Range("I8").Value = 156.2
$ws.Range("K8").Value = 468.6
$ws.Range("M8").Value = -329.6

$ws.Range("H25").Value = 646
$ws.Range("I25").Value = 300
$ws.Range("J25").Value = 992
$ws.Range("K25").Value = 900
$ws.Range("L25").Value = 2976
$ws.Range("M25").Value = -731
$ws.Range("N25").Value = -3314

$ws.Range("H30").Value = 646
$ws.Range("I30").Value = 300
$ws.Range("J30").Value = 992
$ws.Range("K30").Value = 900
$ws.Range("L30").Value = 2976
$ws.Range("M30").Value = -798
$ws.Range("N30").Value = -3180

$ws.Range("H34").Value = 972.7273
$ws.Range("I34").Value = 183.33333
$ws.Range("J34").Value = 1920
$ws.Range("K34").Value = 549.99999
$ws.Range("L34").Value = 5760
$ws.Range("M34").Value = -465.99999
$ws.Range("N34").Value = -5928

$ws.Range("H39").Value = 500
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 1500
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -2088

$ws.Range("H55").Value = 232.8
$ws.Range("I55").Value = 153.5
$ws.Range("K55").Value = 460.5
$ws.Range("M55").Value = -283.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20349.834
$ws.Range("I70").Value = 101900
$ws.Range("J70").Value = 4039.8
$ws.Range("K70").Value = 101900
$ws.Range("L70").Value = 4039.8
$ws.Range("M70").Value = -101630
$ws.Range("N70").Value = -4579.8

$ws.Range("H73").Value = 20349.834
$ws.Range("I73").Value = 101900
$ws.Range("J73").Value = 4039.8
$ws.Range("K73").Value = 101900
$ws.Range("L73").Value = 4039.8
$ws.Range("M73").Value = -100964
$ws.Range("N73").Value = -5911.8

$ws.Range("H80").Value = 19610854
$ws.Range("I80").Value = 66669268
$ws.Range("J80").Value = 3181.5833
$ws.Range("K80").Value = 66669268
$ws.Range("L80").Value = 3181.5833
$ws.Range("M80").Value = -66668270
$ws.Range("N80").Value = -5177.5833

$ws.Range("H83").Value = 19610854
$ws.Range("I83").Value = 66669268
$ws.Range("J83").Value = 3181.5833
$ws.Range("K83").Value = 333346340
$ws.Range("L83").Value = 15907.9165
$ws.Range("M83").Value = -333341348
$ws.Range("N83").Value = -25891.9165

$ws.Range("H97").Value = 1082.1052
$ws.Range("I97").Value = 887.1818
$ws.Range("J97").Value = 1350.125
$ws.Range("K97").Value = 887.1818
$ws.Range("L97").Value = 1350.125
$ws.Range("M97").Value = -391.1818
$ws.Range("N97").Value = -2342.125

$ws.Range("H113").Value = 67890.92999999999
$ws.Range("I113").Value = 143628.72
$ws.Range("K113").Value = 143628.72
$ws.Range("M113").Value = -141458.72

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2364.08
$ws.Range("I96").Value = 1536.6666
$ws.Range("J96").Value = 2625.3684
$ws.Range("K96").Value = 1536.6666
$ws.Range("L96").Value = 2625.3684
$ws.Range("M96").Value = -163.6666
$ws.Range("N96").Value = -5371.368399999999
